$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh: update prices, volume deltas, and row order
# for swapped coins (Cardano/LidoStakedEther, VeChain/Filecoin).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '96.922.30'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +5.11%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.123.91'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.75%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.73'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '608.18'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.66%  '
$ws.Range('E7').Value = '  +2.85%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.386'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('B10').Value = 'LidoStakedEther'
$ws.Range('C10').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.119.24'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.68%  '
$ws.Range('B11').Value = 'Cardano'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.795'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.39%  '
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '96.422.74'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.73%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000240'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.72%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '34.01'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.86%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.34'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.704.06'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.64%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.109.31'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.60%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.58'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.74%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '490.00'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +12.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.47'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.68'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.55%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.0000194'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.45%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.78'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.58%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.53'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.39%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '86.03'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.66'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.79%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.289.77'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.238'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.176'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.44%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.125'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.69%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.844'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -15.52%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '9.05'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.57%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '26.27'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.74%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.38'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -8.39%  '
$ws.Range('E37').Value = '  -3.70%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '491.16'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +5.48%  '
$ws.Range('E39').Value = '  -0.45%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '24.19'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.45%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.439'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E42').Value = '  -2.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.63'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.67%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.21'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.36%  '
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '161.91'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.98%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.701'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.09%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.92'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.00%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '44.32'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.26%  '
$ws.Range('B50').Value = 'Filecoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.35'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.80%  '
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0322'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.00%  '
